$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update watch list symbols
$ws.Range("A2").Value = "OPGN"
$ws.Range("A3").Value = "STOK"
$ws.Range("A4").Value = "CHEK"
$ws.Range("A5").Value = "ABVC"

# Update dates (serial 45377 = 2024-03-26)
$ws.Range("B2").Value = 45377
$ws.Range("B3").Value = 45377
$ws.Range("B4").Value = 45377
$ws.Range("B5").Value = 45377

# Update the selected/active cell shown in the sheet view
$ws.Range("I6").Select()
